$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 8 (pushes /api/rooms... block and everything
# below it down by two rows, and grows the sheet by two rows at the bottom,
# exactly mirroring Excel's native "insert rows" behaviour incl. styles).
$ws.Range("A8:A9").EntireRow.Insert()

# New "/api/signup" endpoint row
$ws.Range("A8").Value = "/api/signup"
$ws.Range("B8").Value = "post/{username,password,email}"
$ws.Range("C8").Value = "{data:{token}}"

# New "/api/signout" endpoint row
$ws.Range("A9").Value = "/api/signout"
$ws.Range("B9").Value = "delete"
$ws.Range("C9").Value = "{data:{message:”success}}"

# Match the author's final cursor position recorded in the workbook view.
$ws.Range("C11").Select()
